$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the placeholder e-mail address used in the sample row.
$ws.Range("A2").Value = "email@example.com"

# Fix the sample "registration street" value (typo/placeholder house number).
$ws.Range("W2").Value = "Хохрякова, 174"

# These two sample cells should use the same font as the other "registration
# address" sample cells (D2 / O2) rather than the one used for the rest of
# the row, so copy just the formatting over (xlPasteFormats = -4122).
$xlPasteFormats = -4122
$ws.Range("D2").Copy()
$ws.Range("A2").PasteSpecial($xlPasteFormats)
$ws.Range("W2").PasteSpecial($xlPasteFormats)
